# The document has three logo pictures living in the header/footer
# stories (wdHeaderFooterFirstPage=2 for the header, wdHeaderFooterPrimary=1
# and wdHeaderFooterFirstPage=2 for the footers). Each InlineShape's
# display name (backed by <wp:docPr name="..."/>) needs to be swapped:
#   - default footer Pearson logo   : image2.png -> image1.png
#   - first-page footer Pearson logo: image2.png -> image1.png
#   - first-page header BTec logo   : image1.jpg -> image2.jpg

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2

# Default (primary) footer - Pearson Edexcel logo (id=1)
$footerPrimary = $sec.Footers.Item($wdHeaderFooterPrimary)
if ($footerPrimary.Exists -and $footerPrimary.Range.InlineShapes.Count -ge 1) {
    $shp = $footerPrimary.Range.InlineShapes.Item(1)
    $shp.Name = "image1.png"
}

# First-page footer - Pearson Edexcel logo (id=2)
$footerFirst = $sec.Footers.Item($wdHeaderFooterFirstPage)
if ($footerFirst.Exists -and $footerFirst.Range.InlineShapes.Count -ge 1) {
    $shp = $footerFirst.Range.InlineShapes.Item(1)
    $shp.Name = "image1.png"
}

# First-page header - BTec logo (id=3)
$headerFirst = $sec.Headers.Item($wdHeaderFooterFirstPage)
if ($headerFirst.Exists -and $headerFirst.Range.InlineShapes.Count -ge 1) {
    $shp = $headerFirst.Range.InlineShapes.Item(1)
    $shp.Name = "image2.jpg"
}
